# Apply the edits described by the commit:
#  - fix G11/I11 to be plain numeric values instead of (mistaken) text entries
#  - add two new compound rows: NO and N2O
#  - update the view (zoom + selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 11 (H2): Cp,b and Cp,c' were stored as text; store as real numbers ---
$ws.Range("G11").Value = 0.00326
$ws.Range("I11").Value = 50200

# --- New row 12: NO (nitrogen (II) oxide), gas ---
# --- New row 13: N2O (nitrogen (II) oxide), gas ---
# (formula/name columns entered first so shared-string order matches: NO, N2O, nitrogen (II) oxide)
$ws.Range("A12").Value = "NO"
$ws.Range("A13").Value = "N2O"
$ws.Range("B12").Value = "nitrogen (II) oxide"
$ws.Range("B13").Value = "nitrogen (II) oxide"

$ws.Range("C12").Value = "g"
$ws.Range("D12").Value = 91.3
$ws.Range("E12").Value = 210.64
$ws.Range("F12").Value = 29.6
$ws.Range("G12").Value = 0.0039
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = -59000
$ws.Range("J12").Value = 298
$ws.Range("K12").Value = 2000

$ws.Range("C13").Value = "g"
$ws.Range("D13").Value = 82
$ws.Range("E13").Value = 219.83
$ws.Range("F13").Value = 45.7
$ws.Range("G13").Value = 0.0086
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = -853000
$ws.Range("J13").Value = 298
$ws.Range("K13").Value = 2000

# --- Update view: zoom to 130% and select I11 ---
$null = $ws.Range("I11").Select()
$excel.ActiveWindow.Zoom = 130
